# Class 2 Weight Update
# Update operating empty weight calculations (column B values shifted from
# the updated column A / WTO calculations).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B16").Value = 1335.1351351351352
$ws.Range("B23").Value = 6.5250839517241932
$ws.Range("B24").Value = 224
$ws.Range("B25").Value = 513.51351351351354
$ws.Range("B26").Value = 41839.193104738762
